$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.979.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.91%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.630.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.18%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'595.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.57%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'155.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.79%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.15%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.95%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  +8.27%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +5.00%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.92%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.26%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'29.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.25%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.0000185"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +21.01%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.092.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.79%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'64.881.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.94%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.636.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.82%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'12.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.22%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.99%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'351.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.16%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'7.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +8.08%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.17%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'68.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.07%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = "'InternetComputer(DFINITY)"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'9.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.77%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "'SuiNetwork"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'1.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.05%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.29%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.164"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.73%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.18%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "'Binance-PegBSC-USD"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.10%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "'PEPE"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0947"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +11.83%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'524.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.99%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +4.32%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +1.79%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +7.98%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'6.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.79%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.424"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.25%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'Monero"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'163.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.88%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'EthereumClassic"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'20.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.78%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.79%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.06%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.05%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'42.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.86%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'164.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.33%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +3.70%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  +4.85%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'23.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.30%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +8.36%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.645"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.77%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0254"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.31%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +2.04%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'19.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.80%  "
$ws.Range("E51").Style = "Normal"
